$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 137; this shifts the existing rows 137:189 down to 138:190
$ws.Rows("137:137").Insert()

# Populate the new row 137 with the new weekly price observation
$ws.Cells.Item(137, 1).Value = 8
$ws.Cells.Item(137, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(137, 3).Value = "Coquimbo"
$ws.Cells.Item(137, 4).Value = 44917
$ws.Cells.Item(137, 5).Value = 4
$ws.Cells.Item(137, 6).Value = 100112044
$ws.Cells.Item(137, 7).Value = "Perejil"
$ws.Cells.Item(137, 8).Value = "Sin especificar"
$ws.Cells.Item(137, 9).Value = "Primera"
$ws.Cells.Item(137, 10).Value = 2000
$ws.Cells.Item(137, 11).Value = 2500
$ws.Cells.Item(137, 12).Value = 3000
$ws.Cells.Item(137, 13).Value = 2750
$ws.Cells.Item(137, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(137, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(137, 16).Value = 1833
$ws.Cells.Item(137, 17).Value = 1.5
$ws.Cells.Item(137, 18).Value = "Hortaliza"

# Make sure the new date cell keeps the same date-time number format as the rest of column D
$ws.Cells.Item(137, 4).NumberFormat = $ws.Cells.Item(138, 4).NumberFormat
